$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix misspelled branch name in cell A6: "Сирдарьинский филиал" -> "Сырдарьинский филиал"
$ws.Range("A6").Value = "Сырдарьинский филиал"

# Update the active cell selection to F24
$ws.Range("F24").Select()
